$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 2 (the "DENISE" row), shifting
# everything else (including the trailing filter-notes row) down by one.
$ws.Rows.Item(2).Insert()

# Fill in the new account row. The leading apostrophe forces the account
# number to be stored as text (preserving the leading zeros) just like the
# other "Conta" cells; ClearFormats() strips the automatic "Text" number
# format Excel would otherwise apply, keeping the cell's style identical to
# its unstyled siblings.
$ws.Cells.Item(2, 1).Value = "'004487140"
$ws.Cells.Item(2, 1).ClearFormats()
$ws.Cells.Item(2, 2).Value = "VALMIR"
$ws.Cells.Item(2, 3).Value = 131883.95
